$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 216.92308
$ws.Range("I28").Value = 147.27272
$ws.Range("K28").Value = 147.27272
$ws.Range("M28").Value = 337.72728
$ws.Range("H32").Value = 2148.6667
$ws.Range("J32").Value = 2094.6924
$ws.Range("L32").Value = 2094.6924
$ws.Range("N32").Value = -2746.6924
$ws.Range("H62").Value = 111116320
$ws.Range("I62").Value = 125002360
$ws.Range("J62").Value = 27980
$ws.Range("K62").Value = 125002360
$ws.Range("L62").Value = 27980
$ws.Range("M62").Value = -125001736
$ws.Range("N62").Value = -29228
$ws.Range("H65").Value = 111116320
$ws.Range("I65").Value = 125002360
$ws.Range("J65").Value = 27980
$ws.Range("K65").Value = 625011800
$ws.Range("L65").Value = 139900
$ws.Range("M65").Value = -625008680
$ws.Range("N65").Value = -146140
$ws.Range("H98").Value = 1401.9615
$ws.Range("J98").Value = 2097.5
$ws.Range("L98").Value = 2097.5
$ws.Range("N98").Value = -5093.5
$ws.Range("H106").Value = 3208.7144
$ws.Range("I106").Value = 2341.6
$ws.Range("J106").Value = 5376.5
$ws.Range("K106").Value = 2341.6
$ws.Range("L106").Value = 5376.5
$ws.Range("M106").Value = -1710.6
$ws.Range("N106").Value = -6638.5
$ws.Range("H107").Value = 1086.909
$ws.Range("I107").Value = 717.3333
$ws.Range("K107").Value = 717.3333
$ws.Range("M107").Value = 1202.6667
$ws.Range("H122").Value = 1401.9615
$ws.Range("J122").Value = 2097.5
$ws.Range("L122").Value = 6292.5
$ws.Range("N122").Value = -11192.5
$ws.Range("H132").Value = 878.61536
$ws.Range("I132").Value = 780
$ws.Range("K132").Value = 2340
$ws.Range("M132").Value = 190
$ws.Range("H138").Value = 2447.322
$ws.Range("I138").Value = 2481.375
$ws.Range("J138").Value = 2406.963
$ws.Range("K138").Value = 7444.125
$ws.Range("L138").Value = 7220.889000000001
$ws.Range("M138").Value = -2304.125
$ws.Range("N138").Value = -17500.889
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 328119.53
$ws.Range("I2").Value = 506304.62
$ws.Range("K2").Value = 506304.62
$ws.Range("M2").Value = -506191.62
$ws.Range("H32").Value = 4462.7163
$ws.Range("I32").Value = 2899.0984
$ws.Range("J32").Value = 20359.5
$ws.Range("K32").Value = 2899.0984
$ws.Range("L32").Value = 20359.5
$ws.Range("M32").Value = -2612.0984
$ws.Range("N32").Value = -20933.5
$ws.Range("H61").Value = 5706.074
$ws.Range("I61").Value = 6198.8335
$ws.Range("K61").Value = 6198.8335
$ws.Range("M61").Value = -5986.8335
$ws.Range("H102").Value = 1384.25
$ws.Range("I102").Value = 1163
$ws.Range("K102").Value = 1163
$ws.Range("M102").Value = 459
$ws.Range("H110").Value = 1031.4615
$ws.Range("I110").Value = 259.9091
$ws.Range("J110").Value = 5275
$ws.Range("K110").Value = 259.9091
$ws.Range("L110").Value = 5275
$ws.Range("M110").Value = 1785.0909
$ws.Range("N110").Value = -9365
$ws.Range("H116").Value = 328119.53
$ws.Range("I116").Value = 506304.62
$ws.Range("K116").Value = 506304.62
$ws.Range("M116").Value = -504010.62
$ws.Range("H136").Value = 5706.074
$ws.Range("I136").Value = 6198.8335
$ws.Range("K136").Value = 18596.5005
$ws.Range("M136").Value = -16046.5005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 328119.53
$ws.Range("I3").Value = 506304.62
$ws.Range("K3").Value = 506304.62
$ws.Range("M3").Value = -506190.62
$ws.Range("H134").Value = 5753.577
$ws.Range("I134").Value = 6403.3335
$ws.Range("K134").Value = 19210.0005
$ws.Range("M134").Value = -16675.0005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H132").Value = 2138.4443
$ws.Range("I132").Value = 1286.5
$ws.Range("K132").Value = 3859.5
$ws.Range("M132").Value = -1329.5
$ws.Range("H134").Value = 1850.4482
$ws.Range("I134").Value = 1521.48
$ws.Range("K134").Value = 4564.440000000001
$ws.Range("M134").Value = -2029.440000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 6088.364
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 6088.364
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 18265.092
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -20345.092
$ws.Range("H140").Value = 3629.5833
$ws.Range("I140").Value = 1211.4
$ws.Range("K140").Value = 3634.2
$ws.Range("M140").Value = 1545.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 20000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H97").Value = 582.4838999999999
$ws.Range("I97").Value = 592.65515
$ws.Range("J97").Value = 435
$ws.Range("K97").Value = 592.65515
$ws.Range("L97").Value = 435
$ws.Range("M97").Value = -96.65515000000005
$ws.Range("N97").Value = -1427
$ws.Range("H122").Value = 2112.5
$ws.Range("I122").Value = 2086.7
$ws.Range("J122").Value = 2177
$ws.Range("K122").Value = 6260.099999999999
$ws.Range("L122").Value = 6531
$ws.Range("M122").Value = -3810.099999999999
$ws.Range("N122").Value = -11431
$ws.Range("H132").Value = 2026838.1
$ws.Range("I132").Value = 2748859.2
$ws.Range("K132").Value = 8246577.600000001
$ws.Range("M132").Value = -8244047.600000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 412142.84
$ws.Range("J2").Value = 90000
$ws.Range("L2").Value = 90000
$ws.Range("N2").Value = -90224
$ws.Range("H22").Value = 2609.6667
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405
$ws.Range("H27").Value = 2609.6667
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593
$ws.Range("H40").Value = 9463.35
$ws.Range("I40").Value = 9198.933999999999
$ws.Range("J40").Value = 10256.6
$ws.Range("K40").Value = 9198.933999999999
$ws.Range("L40").Value = 10256.6
$ws.Range("M40").Value = -9062.933999999999
$ws.Range("N40").Value = -10528.6
$ws.Range("H93").Value = 18519274
$ws.Range("I93").Value = 853.7778
$ws.Range("K93").Value = 853.7778
$ws.Range("M93").Value = 394.2222
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1853.1765
$ws.Range("I136").Value = 1589.6666
$ws.Range("J136").Value = 2149.625
$ws.Range("K136").Value = 4768.9998
$ws.Range("L136").Value = 6448.875
$ws.Range("M136").Value = -2218.9998
$ws.Range("N136").Value = -11548.875
